$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordered list of variable names for rows 2-63 (column A).
# Columns B and C keep whatever value already occupies that row;
# only the variable name (column A) text is updated in place.
$names = @(
    "AAQ_tot",
    "abuse_year",
    "aca_impa",
    "activ_yn",
    "age",
    "anx_mod",
    "anx_sev",
    "assault_emo",
    "assault_phys",
    "assault_sex",
    "audit_tot",
    "belong1",
    "belong2",
    "belong8",
    "belong9",
    "binge_fr",
    "body_sr",
    "BRS_tot",
    "dep_impa",
    "dep_mod",
    "dep_secret",
    "dep_sev",
    "discrim",
    "divers",
    "drugs_yn",
    "dx_adhd",
    "dx_bi",
    "dx_dep",
    "dx_pers",
    "dx_tr",
    "ed_any",
    "env_mh",
    "fincur",
    "finpast",
    "flourish",
    "gad7_impa",
    "gender_noncis",
    "gpa_sr",
    "inf",
    "ins_cover",
    "international",
    "meds_anx",
    "meds_count",
    "meds_dep",
    "meds_mood",
    "meds_sle",
    "meds_sti",
    "military",
    "percneed_cur",
    "persist",
    "psyhx",
    "race",
    "religios",
    "residenc",
    "satisfied_overall",
    "school2_type",
    "sexual",
    "sib_freq",
    "stig_pcv_2",
    "stig_pcv_3",
    "talk",
    "ther_vis"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
}

# Remove the now-obsolete trailing rows (old rows 64-70: sexual, sib_freq,
# stig_pcv_2, stig_pcv_3, talk, ther_vis, wcs_tot) which are no longer
# present after the reshuffle above.
$ws.Range("A64:C70").EntireRow.Delete()

$wb.Save()
